# Changes of 14th June 2022
# Update ShipmentTracking (P), ActualRate (Q) and Result (R) values
# for rows 2-5 on the active worksheet.
#
# ShipmentTracking and ActualRate values look numeric, so the
# NumberFormat is temporarily forced to Text ("@") before assigning the
# value to keep them as literal text (matching the original shared
# string cells), then the style is reset back to the default "Normal"
# style so no new cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("P2") "310109785821"
Set-TextValue $ws.Range("Q2") "$19.10"
$ws.Range("R2").Value = "FAIL"

# Row 3
Set-TextValue $ws.Range("P3") "310109785832"
Set-TextValue $ws.Range("Q3") "$18.06"
$ws.Range("R3").Value = "FAIL"

# Row 4
Set-TextValue $ws.Range("P4") "310109785865"
Set-TextValue $ws.Range("Q4") "$49.94"
$ws.Range("R4").Value = "FAIL"

# Row 5
Set-TextValue $ws.Range("P5") "310109785876"
Set-TextValue $ws.Range("Q5") "$43.56"
$ws.Range("R5").Value = "FAIL"
